$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2024-12-06 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-12-07 Saturday", 2)

# Update the five blocks of three-digit-divided-by-one-digit problems in
# the table. Only the 5 "data" rows (1, 5, 9, 13, 17) contain text; the
# rows between them are blank spacer rows.
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "110÷6=18, 2"
$t.Cell(1,2).Range.Text = "649÷5=129, 4"
$t.Cell(1,3).Range.Text = "540÷5=108, 0"
$t.Cell(1,4).Range.Text = "194÷7=27, 5"
$t.Cell(1,5).Range.Text = "563÷5=112, 3"

$t.Cell(5,1).Range.Text = "223÷7=31, 6"
$t.Cell(5,2).Range.Text = "454÷2=227, 0"
$t.Cell(5,3).Range.Text = "234÷9=26, 0"
$t.Cell(5,4).Range.Text = "407÷9=45, 2"
$t.Cell(5,5).Range.Text = "548÷7=78, 2"

$t.Cell(9,1).Range.Text = "577÷9=64, 1"
$t.Cell(9,2).Range.Text = "643÷2=321, 1"
$t.Cell(9,3).Range.Text = "271÷6=45, 1"
$t.Cell(9,4).Range.Text = "652÷6=108, 4"
$t.Cell(9,5).Range.Text = "245÷5=49, 0"

$t.Cell(13,1).Range.Text = "127÷7=18, 1"
$t.Cell(13,2).Range.Text = "880÷2=440, 0"
$t.Cell(13,3).Range.Text = "410÷8=51, 2"
$t.Cell(13,4).Range.Text = "903÷6=150, 3"
$t.Cell(13,5).Range.Text = "586÷2=293, 0"

$t.Cell(17,1).Range.Text = "803÷2=401, 1"
$t.Cell(17,2).Range.Text = "452÷2=226, 0"
$t.Cell(17,3).Range.Text = "509÷9=56, 5"
$t.Cell(17,4).Range.Text = "864÷8=108, 0"
$t.Cell(17,5).Range.Text = "652÷4=163, 0"
